# Apply the "Commandes SPI" worksheet corrections described in the commit:
#  - add a new column H with "x" markers on the rows that already have a
#    0x.. command code defined (acting as a kind of checklist column)
#  - rename the "0x08" / "Envoie donnee moteur" command to
#    "0x10 + donnee" (the old "Angle moteur" / "***" sub-row is removed)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commandes SPI")
$ws.Activate()

# New column H: mark the rows that already hold a command byte with "x"
$ws.Range("H2").Value  = "x"
$ws.Range("H3").Value  = "x"
$ws.Range("H4").Value  = "x"
$ws.Range("H6").Value  = "x"
$ws.Range("H7").Value  = "x"
$ws.Range("H11").Value = "x"
$ws.Range("H12").Value = "x"
$ws.Range("H13").Value = "x"
$ws.Range("H18").Value = "x"

# "Envoie donnée moteur" now carries the code "0x10 + donnée" instead of
# "0x08", and the old extra "Angle moteur" row underneath it is cleared.
$ws.Range("B18").Value = "0x10 + donnée"
$ws.Range("A19").Value = ""
$ws.Range("B19").Value = ""
$ws.Range("C19").Value = ""

# Match the author's on-screen selection/scroll position after the edit.
$ws.Range("H9").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
